# Auto-update draw results: append the 2025-11-12 "Pick 4" result as row 57.
#
# The sheet stores every value as plain text (dates like "2025-11-12" and
# numeric-looking codes like "251112" must stay text, not become a real
# date/number). A direct `.Value = "..."` assignment lets Excel's
# smart-typing turn those into a date serial / number, so instead we stage
# each value as a formula that evaluates to a text string, copy its
# computed value (not its format) onto the destination cell, then clear the
# scratch cell again. This keeps the new cells' style identical to the
# untouched rows (no extra number format gets attached).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$targetRow = 57
$scratch = $ws.Range("ZZ1000")

function Set-TextValue {
    param($RowNum, $ColLetter, [string]$Text)

    $cellRef = "$ColLetter$RowNum"
    $escaped = $Text -replace '"', '""'

    # Evaluate to a guaranteed string result in the scratch cell ...
    $scratch.Formula = '="' + $escaped + '"'
    # ... then paste only the computed value into the real cell (no format
    # carried over, so the destination keeps its existing/default style) ...
    $scratch.Copy()
    $ws.Range($cellRef).PasteSpecial(-4163)
    # ... and wipe the scratch cell so it leaves no trace in the sheet.
    $scratch.Value = ""
}

Set-TextValue $targetRow "A" "2025-11-12"
Set-TextValue $targetRow "B" "Pick 4"
Set-TextValue $targetRow "C" "251112"
Set-TextValue $targetRow "D" "2-3-0-4"
Set-TextValue $targetRow "E" "2025-11-12T21:40:15.181+04:00"
